$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.876.78'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '3.548.42'
$ws.Range('E3').Value = '  +4.24%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.33'
$ws.Range('E5').Value = '  +3.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.63'
$ws.Range('E6').Value = '  +3.93%  '
$ws.Range('D7').Value = '3.547.48'
$ws.Range('E7').Value = '  +4.26%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +3.43%  '
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.92'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  +4.06%  '
$ws.Range('D13').Value = '4.151.09'
$ws.Range('E13').Value = '  +4.37%  '
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.15'
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').Value = '3.546.66'
$ws.Range('E16').Value = '  +2.94%  '
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '64.793.44'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.04'
$ws.Range('E19').Value = '  +5.68%  '
$ws.Range('E20').Value = '  +6.33%  '
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.40'
$ws.Range('E22').Value = '  +2.80%  '
$ws.Range('E23').Value = '  +6.62%  '
$ws.Range('D24').Value = '3.691.76'
$ws.Range('E24').Value = '  +4.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.27'
$ws.Range('E25').Value = '  +4.20%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +12.80%  '
$ws.Range('E28').Value = '  +7.88%  '
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  +5.35%  '
$ws.Range('E31').Value = '  +4.71%  '
$ws.Range('D32').Value = '3.557.13'
$ws.Range('E32').Value = '  +4.18%  '
$ws.Range('E33').Value = '  +23.57%  '
$ws.Range('E34').Value = '  +4.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.144'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '170.12'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.95'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('E39').Value = '  +7.30%  '
$ws.Range('E40').Value = '  +8.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0807'
$ws.Range('E41').Value = '  +6.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.826'
$ws.Range('E42').Value = '  +3.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.87'
$ws.Range('E43').Value = '  +21.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.64'
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('E46').Value = '  +4.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  +10.69%  '
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('E49').Value = '  +6.73%  '
$ws.Range('D50').Value = '2.448.43'
$ws.Range('E50').Value = '  +12.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('E51').Value = '  +14.80%  '